$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tube")

# Insert a new row at position 5, shifting existing rows down
$ws.Rows.Item(5).Insert()

# Populate the new row 5 with the new monster skill entry
$ws.Range("A5").Value = 1103
$ws.Range("B5").Value = "style_Melee_01"
$ws.Range("C5").Value = "범죄자의 주먹질"
$ws.Range("D5").Value = "style"
$ws.Range("E5").Value = "gangster"
$ws.Range("F5").Value = "D"
$ws.Range("H5").Value = "{(0,5)}"
$ws.Range("I5").Value = "melee"
$ws.Range("K5").Value = 50

# Update the view: zoom to 85%, reset scroll, move selection to A15
$ws.Activate()
$ws.Range("A15").Select() | Out-Null
$excel.ActiveWindow.Zoom = 85
